$d = $word.ActiveDocument

# --- Add new paragraph style "AbstractTitle" (display name "Abstract Title"),
#     inserted conceptually just before the "Abstract" style ---
$title = $d.Styles.Add("AbstractTitle", 1)
$title.NameLocal = "Abstract Title"
$title.BaseStyle = "Normal"
$title.NextParagraphStyle = "Abstract"
$title.QuickStyle = $true

$title.ParagraphFormat.KeepWithNext = $true
$title.ParagraphFormat.KeepTogether = $true
$title.ParagraphFormat.Alignment = 1          # wdAlignParagraphCenter
$title.ParagraphFormat.SpaceBefore = 15       # 300 twentieths-of-a-point
$title.ParagraphFormat.SpaceAfter = 0         # 0 twentieths-of-a-point

$title.Font.Size = 10                         # sz=20 (half-points)
$title.Font.SizeBi = 10                       # szCs=20
$title.Font.Bold = $true
$title.Font.Color = 9067060                   # RGB 345A8A (stored as BGR)

# --- "Abstract" style: reduce space-before from 300 to 100 twentieths (15pt -> 5pt) ---
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- "ImportTok" character style: add green bold color ---
$importTok = $d.Styles("ImportTok")
$importTok.Font.Color = 32768                 # RGB 008000
$importTok.Font.Bold = $true

# --- "BuiltInTok" character style: add green color ---
$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 32768                # RGB 008000
